# "Add files via upload" — the author retyped the unit prices in the
# apiculture sheet as text values ("10 €", "7 €", "5 €") instead of plain
# numbers with currency formatting (entered with a leading apostrophe so
# Excel keeps treating them as text/labels rather than re-parsing them as
# numbers), then left the selection on B7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("apiculture")
$ws.Activate()

# B2:B4 -> "10 €", B5 -> "7 €", B6 -> "5 €" (typed with a leading ' so the
# currency-formatted cells keep their number format but store the value as
# text, i.e. a quoted/quote-prefixed entry).
$ws.Range("B2").Formula = "'10 €"
$ws.Range("B3").Formula = "'10 €"
$ws.Range("B4").Formula = "'10 €"
$ws.Range("B5").Formula = "'7 €"
$ws.Range("B6").Formula = "'5 €"

# Selection ends up on B7 after entering the last value.
$ws.Range("B7").Select()
